$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the variable columns (D, K-T) for data rows 2-19 before overwriting,
# since the edit permutes these values across rows.
$snapshot = @{}
for ($r = 2; $r -le 19; $r++) {
    $row = @{}
    $row["D"] = $ws.Cells.Item($r, 4).Value2
    $row["K"] = $ws.Cells.Item($r, 11).Value2
    $row["L"] = $ws.Cells.Item($r, 12).Value2
    $row["M"] = $ws.Cells.Item($r, 13).Value2
    $row["N"] = $ws.Cells.Item($r, 14).Value2
    $row["O"] = $ws.Cells.Item($r, 15).Value2
    $row["P"] = $ws.Cells.Item($r, 16).Value2
    $row["Q"] = $ws.Cells.Item($r, 17).Value2
    $row["R"] = $ws.Cells.Item($r, 18).Value2
    $row["S"] = $ws.Cells.Item($r, 19).Value2
    $row["T"] = $ws.Cells.Item($r, 20).Value2
    $snapshot[$r] = $row
}

# Row permutation: target row <- source row (derived from the diff)
$mapping = @{}
$mapping[2] = 14
$mapping[3] = 15
$mapping[4] = 5
$mapping[5] = 8
$mapping[6] = 9
$mapping[7] = 2
$mapping[8] = 6
$mapping[9] = 11
$mapping[10] = 12
$mapping[11] = 10
$mapping[12] = 19
$mapping[13] = 7
$mapping[14] = 18
$mapping[15] = 13
$mapping[16] = 3
$mapping[17] = 4
$mapping[18] = 16
$mapping[19] = 17

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $src = $snapshot[$source]
    $ws.Cells.Item($target, 4).Value = $src["D"]
    $ws.Cells.Item($target, 11).Value = $src["K"]
    $ws.Cells.Item($target, 12).Value = $src["L"]
    $ws.Cells.Item($target, 13).Value = $src["M"]
    $ws.Cells.Item($target, 14).Value = $src["N"]
    $ws.Cells.Item($target, 15).Value = $src["O"]
    $ws.Cells.Item($target, 16).Value = $src["P"]
    $ws.Cells.Item($target, 17).Value = $src["Q"]
    $ws.Cells.Item($target, 18).Value = $src["R"]
    $ws.Cells.Item($target, 19).Value = $src["S"]
    $ws.Cells.Item($target, 20).Value = $src["T"]
}
